# Applies the "Fix Sado wording" + Voice QR -> TokiQR edits to the client deck.

$p = $ppt.ActivePresentation

# --- Slide 7 (1-based): Trial Plan description textbox ------------------
# "Voice QR" -> "TokiQR" in the blurb under the Trial Plan heading.
$s7 = $p.Slides.Item(7)
$trialDesc = $s7.Shapes.Item(8)
$trialDesc.TextFrame.TextRange.Text = "Experience three-layer storage with A4 laminate. TokiQR " + [char]8594 + " laminate + NDL deposit + GitHub. Start here."

# --- Slide 8 (1-based): "Relocating to Sado Island" row ------------------
$s8 = $p.Slides.Item(8)

# Shape 19: background pill behind the Sado label - shrink width (right edge
# moves in; left edge untouched). Target EMU 2606040 (205.2pt); the literal
# "205.2" round-trips through the host's single-precision Width setter to
# 205.1999pt/2606039 EMU, so nudge by one float32 ULP to land exactly on
# 2606040 EMU.
$sadoPill = $s8.Shapes.Item(19)
$sadoPill.Width = 205.20001

# Shape 20: the Sado label text box itself - shrink width to match, update
# wording.
$sadoLabel = $s8.Shapes.Item(20)
$sadoLabel.Width = 190.8
$sadoLabel.TextFrame.TextRange.Text = "Establishing base on Sado Island"

# Shape 21: background pill for the next item (Ise Grand Shrine offering) -
# slides left to follow the now-narrower Sado pill; width unchanged. Target
# EMU 3172968 (249.84pt); nudged by one ULP for the same float32 rounding
# reason as above.
$isePill = $s8.Shapes.Item(21)
$isePill.Left = 249.84001

# Shape 22: the Ise label text box - slides left to match; width unchanged.
$iseLabel = $s8.Shapes.Item(22)
$iseLabel.Left = 257.04
